# Auto update ESO load data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated hourly load values for existing rows 8-12 (days 11-15), plus new
# row 13 (day 16) appended to the load data table.
$rows = @{
    8  = @(5101,4879,4787,4779,4759,4919,5376,5875,6299,6467,6401,6299,6201,6138,6055,6052,6077,6166,6357,6326,6101,5812,5597,5257)
    9  = @(4927,4704,4611,4604,4583,4743,5200,5699,6122,6290,6233,6145,6060,6006,5934,5931,5952,6026,6186,6154,5928,5636,5419,5076)
    10 = @(4779,4554,4462,4454,4434,4593,5050,5548,5971,6139,6089,6013,5940,5893,5830,5828,5839,5877,5960,5929,5707,5421,5209,4874)
    11 = @(4577,4307,4196,4183,4176,4205,4360,4721,5215,5548,5633,5600,5565,5522,5470,5452,5472,5521,5567,5478,5230,4940,4811,4575)
    12 = @(4492,4261,4155,4099,4125,4158,4300,4609,5042,5414,5575,5597,5585,5546,5499,5489,5545,5643,5757,5699,5478,5191,4972,4638)
    13 = @(4725,4473,4361,4366,4389,4536,4928,5405,5900,6183,6260,6243,6202,6177,6070,6036,5964,6014,6160,6145,5941,5644,5352,4998)
}

# Row 13 is brand new, so populate the Year/Month/Day columns too.
$ws.Cells.Item(13, 1).Value = 2026
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 16

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 4   # hourly data starts at column D (4)
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
